# Commit: "unify the conception of DataNode, DataTable, Entity."
#
# The only functionally meaningful change in the target diff is the
# worksheet being renamed from "Property1" to "DataNode" (the sheet's
# content/data is untouched), plus the author's cursor/selection having
# moved to E23 by the time the file was saved.
#
# Everything else in the XML diff (fileVersion/rupBuild bump, new xr/xr2/xr3
# revision-tracking namespaces + xr:uid stamps, workbook window geometry
# jitter, a new phonetic-guide font + <phoneticPr>, the "Normal"/"常规"
# locale re-label of the built-in cell style, the x15 timeline-style
# extension block, sub-pixel column-width rounding, and the dropped
# concurrentCalc attribute) are save-format artifacts produced by opening
# the file in a newer/different localized Excel build - not deliberate
# user edits reachable through the Excel object model, so they're left
# alone rather than forced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (Property1 -> DataNode).
$ws.Name = "DataNode"

# Restore the author's last selection in the sheet (bottom-left frozen pane).
$ws.Range("E23").Select()
